# Auto-generated Excel COM-interop script
# Applies per-cell numeric updates (and a few cell clears / one new cell)
# to the Leve profit-calculation sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 4954.8335
$ws.Range("I40").Value = 4496.8887
$ws.Range("J40").Value = 6328.6665
$ws.Range("K40").Value = 4496.8887
$ws.Range("L40").Value = 6328.6665
$ws.Range("M40").Value = -4321.8887
$ws.Range("N40").Value = -6678.6665
# Row 58
$ws.Range("H58").Value = 148.22223
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
# Row 111
$ws.Range("H111").Value = 430.33334
$ws.Range("I111").Value = 430.33334
$ws.Range("K111").Value = 1291.00002
$ws.Range("M111").Value = 1775.99998
# Row 113
$ws.Range("H113").Value = 1495
$ws.Range("I113").Value = 1495
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1495
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1759
$ws.Range("N113").ClearContents()
# Row 116
$ws.Range("H116").Value = 3035.125
$ws.Range("I116").Value = 3056.4
$ws.Range("J116").Value = 2999.6667
$ws.Range("K116").Value = 3056.4
$ws.Range("L116").Value = 2999.6667
$ws.Range("M116").Value = 385.5999999999999
$ws.Range("N116").Value = -9883.6667
# Row 125
$ws.Range("H125").Value = 747.4286
$ws.Range("J125").Value = 825
$ws.Range("L125").Value = 7425
$ws.Range("N125").Value = -12345
# Row 137
$ws.Range("H137").Value = 2191.9697
$ws.Range("I137").Value = 1760.4688
$ws.Range("J137").Value = 16000
$ws.Range("K137").Value = 5281.4064
$ws.Range("L137").Value = 48000
$ws.Range("M137").Value = -2731.4064
$ws.Range("N137").Value = -53100

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 29544.688
$ws.Range("I45").Value = 143358.67
$ws.Range("J45").Value = 3279.923
$ws.Range("K45").Value = 143358.67
$ws.Range("L45").Value = 3279.923
$ws.Range("M45").Value = -142981.67
$ws.Range("N45").Value = -4033.923
# Row 96
$ws.Range("H96").Value = 99998
$ws.Range("J96").Value = 99998
$ws.Range("L96").Value = 99998
$ws.Range("N96").Value = -105490
# Row 97
$ws.Range("H97").Value = 1554.5
$ws.Range("I97").Value = 1611
$ws.Range("K97").Value = 1611
$ws.Range("M97").Value = -1115
# Row 110
$ws.Range("H110").Value = 2077.0667
$ws.Range("I110").Value = 968.5714
$ws.Range("K110").Value = 968.5714
$ws.Range("M110").Value = 1076.4286
# Row 122
$ws.Range("H122").Value = 4184.0312
$ws.Range("I122").Value = 3727.88
$ws.Range("J122").Value = 5813.143
$ws.Range("K122").Value = 11183.64
$ws.Range("L122").Value = 17439.429
$ws.Range("M122").Value = -8733.639999999999
$ws.Range("N122").Value = -22339.429
# Row 132
$ws.Range("H132").Value = 1811.2727
$ws.Range("I132").Value = 1416.2
$ws.Range("K132").Value = 4248.6
$ws.Range("M132").Value = -1718.6

$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 52345.25
$ws.Range("J82").Value = 91921.664
$ws.Range("L82").Value = 91921.664
$ws.Range("N82").Value = -92687.664
# Row 85
$ws.Range("H85").Value = 52345.25
$ws.Range("J85").Value = 91921.664
$ws.Range("L85").Value = 91921.664
$ws.Range("N85").Value = -94573.664
# Row 86
$ws.Range("H86").Value = 5613.625
$ws.Range("I86").Value = 5429.4546
$ws.Range("K86").Value = 5429.4546
$ws.Range("M86").Value = -4306.4546
# Row 89
$ws.Range("H89").Value = 5613.625
$ws.Range("I89").Value = 5429.4546
$ws.Range("K89").Value = 27147.273
$ws.Range("M89").Value = -21531.273
# Row 134
$ws.Range("H134").Value = 2163.611
$ws.Range("I134").Value = 1613.0333
$ws.Range("K134").Value = 4839.0999
$ws.Range("M134").Value = -2304.0999

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 2699.8
$ws.Range("I22").Value = 1299.5
$ws.Range("K22").Value = 1299.5
$ws.Range("M22").Value = -949.5
# Row 31
$ws.Range("H31").Value = 3591.75
$ws.Range("I31").Value = 2808.392
$ws.Range("J31").Value = 5494.1904
$ws.Range("K31").Value = 2808.392
$ws.Range("L31").Value = 5494.1904
$ws.Range("M31").Value = -2513.392
$ws.Range("N31").Value = -6084.1904
# Row 34
$ws.Range("H34").Value = 3591.75
$ws.Range("I34").Value = 2808.392
$ws.Range("J34").Value = 5494.1904
$ws.Range("K34").Value = 2808.392
$ws.Range("L34").Value = 5494.1904
$ws.Range("M34").Value = -2606.392
$ws.Range("N34").Value = -5898.1904
# Row 52
$ws.Range("H52").Value = 95836.5
$ws.Range("J52").Value = 95836.5
$ws.Range("L52").Value = 95836.5
$ws.Range("N52").Value = -96424.5
# Row 99
$ws.Range("H99").Value = 6404.6
$ws.Range("I99").Value = 6404.6
$ws.Range("K99").Value = 6404.6
$ws.Range("M99").Value = -4906.6
# Row 126
$ws.Range("H126").Value = 6404.6
$ws.Range("I126").Value = 6404.6
$ws.Range("K126").Value = 19213.8
$ws.Range("M126").Value = -16743.8
# Row 141
$ws.Range("H141").Value = 274572.2
$ws.Range("J141").Value = 274572.2
$ws.Range("L141").Value = 274572.2
$ws.Range("N141").Value = -284932.2

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 2362
$ws.Range("J122").Value = 2452.625
$ws.Range("L122").Value = 22073.625
$ws.Range("N122").Value = -26973.625
# Row 132
$ws.Range("H132").Value = 3494.68
$ws.Range("I132").Value = 2759.6
$ws.Range("J132").Value = 3678.45
$ws.Range("K132").Value = 24836.4
$ws.Range("L132").Value = 33106.05
$ws.Range("M132").Value = -22306.4
$ws.Range("N132").Value = -38166.05
# Row 137
$ws.Range("H137").Value = 4174.579
$ws.Range("I137").Value = 3752.5
$ws.Range("J137").Value = 5356.4
$ws.Range("K137").Value = 11257.5
$ws.Range("L137").Value = 16069.2
$ws.Range("M137").Value = -6157.5
$ws.Range("N137").Value = -26269.2

$ws = $wb.Worksheets.Item("GSM")
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
# Row 113
$ws.Range("H113").Value = 3949
$ws.Range("I113").Value = 3767.8708
$ws.Range("K113").Value = 3767.8708
$ws.Range("M113").Value = -1597.8708
# Row 122
$ws.Range("H122").Value = 2087.4
$ws.Range("I122").Value = 2079.5
$ws.Range("K122").Value = 6238.5
$ws.Range("M122").Value = -3788.5
# Row 132
$ws.Range("H132").Value = 2724.5593
$ws.Range("I132").Value = 2249.9148
$ws.Range("K132").Value = 6749.7444
$ws.Range("M132").Value = -4219.7444

$ws = $wb.Worksheets.Item("LTW")
# Row 12
$ws.Range("H12").Value = 25000
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
# Row 22
$ws.Range("H22").Value = 2026.5
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
# Row 27
$ws.Range("H27").Value = 2026.5
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
# Row 46
$ws.Range("H46").Value = 1523.4517
$ws.Range("J46").Value = 1382.8518
$ws.Range("L46").Value = 1382.8518
$ws.Range("N46").Value = -1758.8518
# Row 122
$ws.Range("H122").Value = 5507.077
$ws.Range("I122").Value = 5528.3
$ws.Range("K122").Value = 16584.9
$ws.Range("M122").Value = -14134.9

$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value = 46996
$ws.Range("J70").Value = 46996
$ws.Range("L70").Value = 46996
$ws.Range("N70").Value = -47626
# Row 73
$ws.Range("H73").Value = 46996
$ws.Range("J73").Value = 46996
$ws.Range("L73").Value = 46996
$ws.Range("N73").Value = -49180
# Row 122
$ws.Range("H122").Value = 17860842
$ws.Range("I122").Value = 3429
$ws.Range("K122").Value = 10287
$ws.Range("M122").Value = -7837

